# Apply the ValidLogin changes to the workbook.
#
# 1. Rename the single worksheet from "test1" to "ValidLogin".
# 2. Populate a small username/password table used by the new
#    ValidLogin.java test (headers in row 1, sample credentials in row 2).
# 3. Leave the active selection on C17, matching the author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "test1" -> "ValidLogin"
$ws.Name = "ValidLogin"

# 2. Write header row + credential row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# 3. Restore the selection the author left the sheet on
$ws.Range("C17").Select()
